# Add team record (Wins / Losses / Ties) columns to the PIT_1999 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF.
# Copy the formatting from the existing header cell (AC1) so the new
# headers match the bold/bordered/centered look of the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-47 all share the same constant team record values
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 78  # AD
    $ws.Cells.Item($r, 31).Value = 83  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
